$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "Finished?" markers in column I for the sounds/animation rows
$ws.Range("I18").Value = "Yes"
$ws.Range("I22").Value = "WIP"
$ws.Range("I24").Value = "Yes"

# Update existing "Finished?" status in column E
$ws.Range("E30").Value = "Yes"
$ws.Range("E36").Value = "WIP"
$ws.Range("E40").Value = "Yes"

# Update the view state: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E31").Select()
